$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Strip the leading space from the header labels in B1 and D1 (" A" -> "A", " C" -> "C")
$ws.Range("B1").Value = "A"
$ws.Range("D1").Value = "C"

# Give the boolean cells D2/D4 explicit TRUE()/FALSE() formulas (cached boolean values unchanged)
$ws.Range("D2").Formula = "=TRUE()"
$ws.Range("D4").Formula = "=FALSE()"

# Move the active selection/cursor to D4 (last cell touched)
[void]$ws.Range("D4").Select()
